$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.178.68"
$ws.Range("D3").Value = "1.660.36"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'216.93"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").Value = "'0.5146"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").Value = "'0.06270"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").Value = "'0.07741"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "1.664.50"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'4.447"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "1.887.54"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "'0.5437"
$ws.Range("D16").Value = "0.0₅8098"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "'64.89"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "26.195.23"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'4.630"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "'192.06"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "'10.08"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "'6.018"
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'140.05"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "'0.1223"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").Value = "'7.226"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'16.15"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "'0.05968"
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("D31").Value = "'1.272"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "'3.569"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("D35").Value = "'0.9659"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").Value = "'2.422"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'0.5667"
$ws.Range("E38").Value = "  -7.96%  "
$ws.Range("D39").Value = "'0.01589"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "'5.968"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "'0.8571"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "1.014.99"
$ws.Range("E43").Value = "  -7.13%  "
$ws.Range("D44").Value = "'100.24"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "1.802.08"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "'56.64"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'7.993"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -3.87%  "

Write-Host "Applied cryptos update"
